$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.002.99"
$ws.Range("E2").Value = "  -0.65%  "

$ws.Range("D3").Value = "2.354.32"
$ws.Range("E3").Value = "  -0.36%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.680"
$ws.Range("E5").Value = "  +0.39%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "239.34"
$ws.Range("E6").Value = "  +0.16%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.13"
$ws.Range("E7").Value = "  +0.87%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.598"
$ws.Range("E9").Value = "  +9.37%  "

$ws.Range("E10").Value = "  -0.26%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.31"
$ws.Range("E11").Value = "  +0.18%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "32.55"
$ws.Range("E12").Value = "  +10.43%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.31"
$ws.Range("E13").Value = "  +8.75%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.107"
$ws.Range("E14").Value = "  -0.11%  "

$ws.Range("D15").Value = "2.702.00"
$ws.Range("E15").Value = "  -0.45%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "16.60"
$ws.Range("E16").Value = "  -1.48%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.899"
$ws.Range("E17").Value = "  -0.58%  "

$ws.Range("D18").Value = "2.358.36"
$ws.Range("E18").Value = "  -0.06%  "

$ws.Range("D19").Value = "43.891.88"
$ws.Range("E19").Value = "  -0.63%  "

$ws.Range("E20").Value = "  -0.93%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.75"
$ws.Range("E21").Value = "  +4.55%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "76.90"
$ws.Range("E22").Value = "  -1.55%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "258.52"
$ws.Range("E23").Value = "  +0.89%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.95"
$ws.Range("E24").Value = "  +22.22%  "

$ws.Range("E25").Value = "  +0.01%  "

$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.49"
$ws.Range("E26").Value = "  -1.37%  "

$ws.Range("B27").Value = "WEMIXToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.65"
$ws.Range("E27").Value = "  -2.89%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.75"
$ws.Range("E28").Value = "  +2.08%  "

$ws.Range("E29").Value = "  -0.51%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.64"
$ws.Range("E30").Value = "  +0.62%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "175.56"
$ws.Range("E31").Value = "  +1.44%  "

$ws.Range("E32").Value = "  -2.57%  "

$ws.Range("E33").Value = "  +2.61%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0764"
$ws.Range("E34").Value = "  +3.66%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.22"
$ws.Range("E35").Value = "  +0.08%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.46"
$ws.Range("E36").Value = "  +4.61%  "

$ws.Range("E37").Value = "  -4.49%  "

$ws.Range("E38").Value = "  -3.82%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.30"
$ws.Range("E39").Value = "  -2.96%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0278"
$ws.Range("E40").Value = "  +2.46%  "

$ws.Range("E41").Value = "  +12.52%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.207"
$ws.Range("E42").Value = "  +13.31%  "

$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "18.94"
$ws.Range("E43").Value = "  -3.79%  "

$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.01"
$ws.Range("E44").Value = "  +1.41%  "

$ws.Range("E45").Value = "  -0.05%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.72"
$ws.Range("E46").Value = "  +4.99%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.52"
$ws.Range("E47").Value = "  +6.81%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "58.12"
$ws.Range("E48").Value = "  +9.72%  "

$ws.Range("E49").Value = "  -0.45%  "

$ws.Range("E50").Value = "  +0.09%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "100.03"
$ws.Range("E51").Value = "  +1.29%  "
